$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E ("Preconditions") gains "None" for the first several test cases ---
$ws.Range("E7").Value  = "None"
$ws.Range("E8").Value  = "None"
$ws.Range("E9").Value  = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "None"

# --- Row 11 (Test case 5): rewrite Method Inputs / Expected Result ---
$ws.Range("F11").Value = 'client= Client(123,"Lily ", "Green", "lilygreengmail")'
$ws.Range("G11").Value = 'client.email_address = "email@pixell-river.com" '

# --- Row 12 (Test case 6: client_number getter) ---
$ws.Range("E12").Value = "client_number=123"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "self.client.client_number = 123"

# --- Row 13 (Test case 7: first_name getter) ---
$ws.Range("E13").Value = 'first_name="Lily"'
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = 'self.client.first_name = "Lily"'

# --- Row 14 (Test case 8: last_name getter) ---
$ws.Range("E14").Value = 'last_name="Green"'
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = 'self.client.last_name = "Green"'

# --- Row 15 (Test case 9: email_address getter) ---
$ws.Range("E15").Value = 'email_address="lilygreen@gmail.com"'
$ws.Range("F15").Value = "None"
$ws.Range("G15").Value = 'self.client.email_address="lilygreen@gmail.com"'

# --- Row 16 (Test case 10: __str__) ---
$ws.Range("E16").Value = 'client_number=123,             first_name="Lily",      last_name="Green",        email_address="lilygreen@gmail.com"'
$ws.Range("F16").Value = "None"
$ws.Range("G16").Value = '"Green, Lily[123]-lilygreen@gmail.com"'

# E16/F16 pick up the same bold "table body" font used by the rest of column E-G (style 7)
$ws.Range("E16").Font.Bold = $true
$ws.Range("F16").Font.Bold = $true

# --- Row height tweaks that came along with the re-wrapped text ---
$ws.Rows.Item(11).RowHeight = 54.75
$ws.Rows.Item(12).RowHeight = 60.4

# --- Column width tweaks for E:G ---
# (ColumnWidth set here + an internal ~5/6-character padding is what ends up
#  stored as the OOXML <col width>; values below are chosen so the saved
#  width lands as close as possible to 20.46484375 / 22.86328125 / 33.265625)
$ws.Columns.Item(5).ColumnWidth = 19.666666666666668
$ws.Columns.Item(6).ColumnWidth = 22
$ws.Columns.Item(7).ColumnWidth = 32.5

# --- Selection / scroll position left by the author when they saved ---
$ws.Range("F16").Select()
